# Update data collection model
#
# DataCollectionGroup / DataCollectionGroup1 (sheet2.xml / sheet8.xml):
#   Replace columns G:H (collectionigStartVersion, collectionigEndVersion)
#   with two new columns inserted before "domain" (E:F):
#   standardStartVersion, standardEndVersion.
#
# DataCollectionItem / DataCollectionItem1 (sheet3.xml / sheet9.xml):
#   Remove column Q (cdashigCore); sdtmTarget shifts from R to Q.

$wb = $excel.ActiveWorkbook

$groupSheets = @("DataCollectionGroup", "DataCollectionGroup1")
foreach ($name in $groupSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns("G:H").Delete()
    $ws.Columns("E:F").Insert()
    $ws.Range("E1").Value = "standardStartVersion"
    $ws.Range("F1").Value = "standardEndVersion"
}

$itemSheets = @("DataCollectionItem", "DataCollectionItem1")
foreach ($name in $itemSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns("Q").Delete()
}
